# Sprint Backlog Week2 update
# 1) Reassign tasks in rows 16-20 (Sprint item "Implement remove/add/edit..." rows)
#    from "TBD" to "Vitor".
# 2) Fill in Week 2 progress column (F) for rows 12-20 with this week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Re-assign owner for rows 16-20 from TBD -> Vitor
$ws.Range("B16:B20").Value = "Vitor"

# 2) Week 2 (column F) task-completion values for rows 12-20
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = 1
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 1

# Leave the selection where the author ended up after editing F20
$ws.Range("F20").Select()
